$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 7 new rows before row 9 (pushes old row9 "totals" -> row16, old row10 "footer" -> row17)
$ws.Rows.Item(9).Resize(7).Insert(-4121)

# 2) Clone the formatting of row 8 (same layout/styles as every item row) onto the 7 new rows
$ws.Range("A8:Q8").Copy()
for ($r = 9; $r -le 15; $r++) {
    $ws.Range("A" + $r + ":Q" + $r).PasteSpecial(-4122)
}

# 3) Row heights for the new rows (matches the per-row custom heights from the template)
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75

# 4) Re-create the per-row merges (A:B, C:G, H:K, L:M, N:O) for each of the new rows
for ($r = 9; $r -le 15; $r++) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

Write-Host "structure done"

# 5) Item-row data (A = running number, C = item name, H = "out-of-stock:on-order"
#    counter, L = reorder-level flag, N = unit price, P = sale price, Q = transaction count)
$items = @(
    @{ Row = 7;  A = 1; C = "  1+1 INFINITY CARE 50+ SPF LOTION"; H = "0:0"; L = "0"; N = "264.00"; P = "264.0000"; Q = "1:0" },
    @{ Row = 8;  A = 2; C = "ANTI-COX II 15MG/3ML 6 AMP";          H = "0:1"; L = "1"; N = "78.00";  P = "12.4800"; Q = "0:1" },
    @{ Row = 9;  A = 3; C = "DIAMICRON 60MG M.R. 30 SCORED TAB";   H = "3:1"; L = "1"; N = "108.00"; P = "108.0000"; Q = "1:0" },
    @{ Row = 10; A = 4; C = "EMPACOZA TRIO XR 25/5/1000  30TAB";   H = "1:1"; L = "0"; N = "396.00"; P = "130.6800"; Q = "0:1" },
    @{ Row = 11; A = 5; C = "ERASTAPEX 20 MG 30 F.C.TAB.";         H = "1:0"; L = "1"; N = "75.00";  P = "24.7500"; Q = "0:1" },
    @{ Row = 12; A = 6; C = "METACARDIA MR 35MG 30 F.C. TAB.";     H = "0:2"; L = "1"; N = "60.00";  P = "60.0000"; Q = "1:0" },
    @{ Row = 13; A = 7; C = "TERRAMYCIN EYE OINT. 5 GM";           H = "2:0"; L = "1"; N = "28.00";  P = "28.0000"; Q = "1:0" },
    @{ Row = 14; A = 8; C = "ZYROVAZET 10/20MG 30 F.C. TABLETS";   H = "0:2"; L = "1"; N = "294.00"; P = "97.0200"; Q = "0:1" },
    @{ Row = 15; A = 9; C = "سرنجات 3 سم";                         H = "0:0"; L = "0"; N = "2.00";   P = "2.0000";  Q = "1:0" }
)

foreach ($it in $items) {
    $r = $it.Row
    $ws.Range("A" + $r).Value = $it.A
    $ws.Range("C" + $r).Value = $it.C
    $ws.Range("H" + $r).Value = $it.H
    $ws.Range("L" + $r).Value = $it.L
    $ws.Range("N" + $r).Value = $it.N
    $ws.Range("P" + $r).Value = $it.P
    $ws.Range("Q" + $r).Value = $it.Q
}

# 6) Totals row (now row 16) - literal sum of the "sale price" column carried over
#    verbatim from the source export (keeps the exact floating point representation)
$ws.Range("P16").Value = 726.92999999999995

# 7) Footer row (now row 17) - refreshed export timestamp, page count & signature
$ws.Range("A17").Value = "Monday, 21 July, 2025 9:56 AM"
$ws.Range("G17").Value = "1/1"
$ws.Range("K17").Value = "developed by : Abdelaziz Talaat"

Write-Host "data done"
